$d = $word.ActiveDocument

# 1. Drop the old "_GoBack" bookmark left over near "Контроль ... занятости компьютеров".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Remove the leading sentence "Должен быть хотя бы один человек: администратор. "
#    from section 3.3, keeping "Требования к администратору: ..." intact.
$d.Content.Find.Execute(
    "Должен быть хотя бы один человек: администратор. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    2)

# 3. Word re-drops "_GoBack" at the location of the latest edit -- recreate it right
#    before "Требования к администратору: " to match.
$rng = $d.Content
$rng.Find.Execute(
    "Требования к администратору: ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($rng.Start, $rng.Start)
$d.Bookmarks.Add("_GoBack", $insertPoint)
